$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# --- Column H width (bestFit grew to match new "-100" values) ---
$ws.Columns.Item(8).ColumnWidth = 6.65

# --- Data table updates ---
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4163)
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("H14").Value = -100
$ws.Range("C15").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 23
$ws.Range("K15").Value = 27.777777777777
$ws.Range("L15").Value = 155.555555555556
$ws.Range("M15").Value = 666.666666666667
$ws.Range("N15").Value = 27.777777777777
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = -46.153846153846
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 42
$ws.Range("H16").Value = -42.857142857142
$ws.Range("I16").Value = 179
$ws.Range("J16").Value = 204
$ws.Range("K16").Value = -12.254901960784
$ws.Range("L16").Value = -29.803921568627
$ws.Range("M16").Value = 155.714285714286
$ws.Range("N16").Value = -86.31498470948
$ws.Range("C17").Value = 10
$ws.Range("E17").Value = -23.076923076923
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 289
$ws.Range("J17").Value = 279
$ws.Range("K17").Value = 3.584229390681
$ws.Range("L17").Value = 18.442622950819
$ws.Range("M17").Value = 197.938144329897
$ws.Range("N17").Value = -11.349693251533
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = 30
$ws.Range("I18").Value = 212
$ws.Range("J18").Value = 193
$ws.Range("K18").Value = 9.844559585492
$ws.Range("L18").Value = -4.932735426008
$ws.Range("M18").Value = 20.454545454545
$ws.Range("N18").Value = -84.143605086013
$ws.Range("C19").Value = 37
$ws.Range("D19").Value = 45
$ws.Range("E19").Value = -17.777777777777
$ws.Range("F19").Value = 135
$ws.Range("G19").Value = 143
$ws.Range("H19").Value = -5.594405594405
$ws.Range("I19").Value = 906
$ws.Range("J19").Value = 1048
$ws.Range("K19").Value = -13.54961832061
$ws.Range("L19").Value = -27.635782747603
$ws.Range("M19").Value = -21.829163071613
$ws.Range("N19").Value = -81.038091251569
$ws.Range("D20").Value = 2
$ws.Range("C14").Copy()
$ws.Range("F20").PasteSpecial(-4163)
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -100
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = -57.142857142857
$ws.Range("L20").Value = -61.290322580645
$ws.Range("M20").Value = -14.285714285714
$ws.Range("N20").Value = -93.513513513513
$ws.Range("C21").Value = 58
$ws.Range("D21").Value = 76
$ws.Range("E21").Value = -23.684210526315
$ws.Range("G21").Value = 250
$ws.Range("H21").Value = -9.6
$ws.Range("I21").Value = 1622
$ws.Range("J21").Value = 1773
$ws.Range("K21").Value = -8.516638465877
$ws.Range("L21").Value = -19.54365079365
$ws.Range("M21").Value = 6.780776826859
$ws.Range("N21").Value = -79.620555346149
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -20
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 18.181818181818
$ws.Range("I22").Value = 111
$ws.Range("J22").Value = 94
$ws.Range("K22").Value = 18.085106382978
$ws.Range("L22").Value = -6.72268907563
$ws.Range("M22").Value = 65.671641791044
$ws.Range("C24").Value = 66
$ws.Range("D24").Value = 85
$ws.Range("E24").Value = -22.35294117647
$ws.Range("F24").Value = 299
$ws.Range("G24").Value = 343
$ws.Range("H24").Value = -12.827988338192
$ws.Range("I24").Value = 2059
$ws.Range("J24").Value = 2225
$ws.Range("K24").Value = -7.460674157303
$ws.Range("L24").Value = 0.881920627143
$ws.Range("M24").Value = -14.029227557411
$ws.Range("D25").Value = 77
$ws.Range("E25").Value = -23.376623376623
$ws.Range("F25").Value = 238
$ws.Range("G25").Value = 314
$ws.Range("H25").Value = -24.203821656051
$ws.Range("I25").Value = 1768
$ws.Range("J25").Value = 1970
$ws.Range("K25").Value = -10.253807106599
$ws.Range("L25").Value = -7.579717720857
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 82
$ws.Range("G26").Value = 85
$ws.Range("H26").Value = -3.529411764705
$ws.Range("I26").Value = 560
$ws.Range("J26").Value = 529
$ws.Range("K26").Value = 5.86011342155
$ws.Range("L26").Value = 5.66037735849
$ws.Range("M26").Value = 86.666666666666
$ws.Range("C27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("H14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 25
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = 8.695652173913
$ws.Range("L27").Value = 56.25
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 225
$ws.Range("F28").Value = 38
$ws.Range("G28").Value = 22
$ws.Range("H28").Value = 72.727272727272
$ws.Range("I28").Value = 142
$ws.Range("J28").Value = 111
$ws.Range("K28").Value = 27.927927927927
$ws.Range("L28").Value = 23.478260869565
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4163)
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4163)
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100
$ws.Range("D31").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$ws.Range("H14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = -66.666666666666
$ws.Range("J31").Value = 12
$ws.Range("K31").Value = -25
$ws.Range("L31").Value = -10
